$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the D, N, O, P, Q, R, S, T column values between row 2 and row 3
# Row 2 new values (previously in row 3)
$ws.Range("D2").Value = 44330
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 15500
$ws.Range("Q2").Value = "`$/caja 18 kilos granel"
$ws.Range("R2").Value = "Provincia de Curicó"
$ws.Range("S2").Value = 861
$ws.Range("T2").Value = 18

# Row 3 new values (previously in row 2)
$ws.Range("D3").Value = 44334
$ws.Range("N3").Value = 11000
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 11500
$ws.Range("Q3").Value = "`$/caja 12 kilos granel"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 11500
$ws.Range("T3").Value = 1
